$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17 (pushes existing rows 17..83 down to 18..84,
# and extends the sheet dimension to A1:R84 automatically).
$ws.Rows.Item(17).Insert()

# Populate the newly-inserted row 17 with a new weekly record (same market /
# product as the row that used to be here, but a new sample date).
$ws.Range("A17").Value = 7
$ws.Range("B17").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C17").Value = "Ñuble"
$ws.Range("D17").Value = 45189
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = 100112026
$ws.Range("G17").Value = "Haba"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 15000
$ws.Range("N17").Value = "$/saco 25 kilos"
$ws.Range("O17").Value = "Provincia de Diguillín"
$ws.Range("P17").Value = 600
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = "Hortaliza"
